$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('AT2').Value = 'When miss emily grierson died, our whole town went to her funeral: the men through a sort of respectful affection for a fallen monument, the women mostly out of curiosity to see the inside of her house, which no one save an old man-servant--a combined gardener and cook--had seen in at least ten years. it was a big, squarish frame house that had once been white, decorated with cupolas and spires and scrolled balconies in the heavily lightsome style of the seventies, set on what had once been our most select street. but garages and cotton gins had encroached and obliterated even the august names of that neighborhood; only miss emily''s house was left, lifting its stubborn and coquettish decay above the cotton wagons and the gasoline pumps-an eyesore among eyesores. and now miss emily had gone to join the representatives of those august names where they lay in the cedar-bemused cemetery among '
$ws.Range('AT3').Value = 'the ranked and anonymous graves of union and confederate soldiers who fell at the battle of jefferson. alive, miss emily had been a tradition, a duty, and a care; a sort of hereditary obligation upon the town, dating from '
$ws.Range('AT4').Value = 'that day in 1894 when colonel sartoris, the mayor--he who fathered the edict that no negro woman should appear on the streets without an apron-remitted her taxes, '
$ws.Range('AT5').Value = 'the dispensation dating from the death of her father on into perpetuity. not that miss emily would have accepted charity. colonel sartoris invented an involved tale to the effect that miss emily''s father had loaned money to the town, which the town, as a matter of business, preferred this way of repaying. only a man of colonel sartoris'' generation and thought could have invented it, and only a woman could have believed it. '
$ws.Range('AT6').Value = 'When the next generation, with its more modern ideas, became mayors and aldermen, this arrangement created some little dissatisfaction. on the first of the year they mailed her a tax notice. february came, and there was no reply. they wrote her a formal letter, asking her to call at the sheriff''s office at her convenience. a week later the mayor wrote her himself, offering to call or to send his car for her, and received in reply a note on paper of an archaic shape, in a thin, flowing calligraphy in faded ink, to the effect that she no longer went out at all. the tax notice was also enclosed, without comment.  they called a special meeting of the board of aldermen. '
$ws.Range('AT7').Value = 'A deputation waited upon her, knocked at the door through which no visitor had passed since she ceased giving china-painting lessons eight or ten years earlier. they were admitted by the old negro into a dim hall from which a stairway mounted into still more shadow. it smelled of dust and disuse-a close, dank smell. the negro led them into the parlor. it was furnished in heavy, leather-covered furniture. when the negro opened the blinds of one window, they could see that the leather was cracked; and when they sat down, a faint dust rose sluggishly about their thighs, spinning with slow motes in the single sun-ray. on a tarnished gilt easel before the fireplace stood a crayon portrait of miss emily''s father. they rose when she entered--a small, fat woman in black, with a thin gold chain descending to her waist and vanishing into her belt, leaning on an ebony cane with a tarnished gold head. her skeleton was small and spare; perhaps that was why what would have been merely plumpness in another was obesity in her. she looked bloated, like a body long submerged in motionless water, and of that pallid hue. her eyes, lost in the fatty ridges of her face, looked like two small pieces of coal pressed into a lump of dough as they moved from one face to another while the visitors stated their errand. she did not ask them to sit. she just stood in the door and listened quietly until the spokesman came to a stumbling halt. then they could hear the invisible watch ticking at the end of the gold chain. her voice was dry and cold. "i have no taxes in jefferson. colonel sartoris explained it to me. perhaps one of you can gain access to the city records and satisfy yourselves." "but we have. we are the city authorities, miss emily. didn''t you get a notice from the sheriff, signed by him?" "i received a paper, yes," miss emily said. "perhaps he considers himself the sheriff . . . i have no taxes in jefferson." "but there is nothing on the books to show that, you see we must go by the--" "see colonel sartoris. i have no taxes in jefferson." "but, miss emily--"  "see colonel sartoris." (colonel sartoris had been dead almost ten years.) "i have no taxes in jefferson. tobe!" the negro appeared. "show these gentlemen out."  ii'
$ws.Range('AT8').Value = ' so she vanquished them, horse and foot, just as she had vanquished their fathers thirty years before about the smell. that was two years after her father''s death and a short time after her sweetheart--the one we believed would marry her --had deserted her. after her father''s death she went out very little; after her sweetheart went away, people hardly saw her at all. a few of the ladies had the temerity to call, but were not received, and the only sign of life about the place was the negro man--a young man then--going in and out with a market basket. "just as if a man--any man--could keep a kitchen properly, "the ladies said; so they were not surprised when the smell developed. it was another link between the gross, teeming world and the high and mighty griersons. '
$ws.Range('AT9').Value = 'A neighbor, a woman, complained to the mayor, judge stevens, eighty years old. "but what will you have me do about it, madam?" he said. "why, send her word to stop it," the woman said. "isn''t there a law? " "i''m sure that won''t be necessary," judge stevens said. "it''s probably just a snake or a rat that nigger of hers killed in the yard. i''ll speak to him about it." the next day he received two more complaints, one from a man who came in diffident deprecation. "we really must do something about it, judge. i''d be the last one in the world to bother miss emily, but we''ve got to do something." that night the board of aldermen met--three graybeards and one younger man, a member of the rising generation. "it''s simple enough," he said. "send her word to have her place cleaned up. give her a certain time to do it in, and if she don''t. .." "dammit, sir," judge stevens said, "will you accuse a lady to her face of smelling bad?" '
$ws.Range('AT10').Value = ' so the next night, after midnight, four men crossed miss emily''s lawn and slunk about the house like burglars, sniffing along the base of the brickwork and at the cellar openings while one of them performed a regular sowing motion with his hand out of a sack slung from his shoulder. they broke open the cellar door and sprinkled lime there, and in all the outbuildings. as they recrossed the lawn, a window that had been dark was lighted and miss emily sat in it, the light behind her, and her upright torso motionless as that of an idol. they crept quietly across the lawn and into the shadow of the locusts that lined the street. after a week or two the smell went away. that was when people had begun to feel really sorry for her. '
$ws.Range('AT11').Value = 'People in our town, remembering how old lady wyatt, her great-aunt, had gone completely crazy at last, believed that the griersons held themselves a little too high for what they really were. '
$ws.Range('AT12').Value = 'None of the young men were quite good enough for miss emily and such. we had long thought of them as a tableau, miss emily a slender figure in white in the background, her father a spraddled silhouette in the foreground, his back to her and clutching a horsewhip, the two of them framed by the back-flung front door. so when she got to be thirty and was still single, we were not pleased exactly, but vindicated; even with insanity in the family she wouldn''t have turned down all of her chances if they had really materialized. '
$ws.Range('AT13').Value = 'When her father died, it got about that the house was all that was left to her; and in a way, people were glad. at last they could pity miss emily. being left alone, and a pauper, she had become humanized. now she too would know the old thrill and the old despair of a penny more or less. the day after his death all the ladies prepared to call at the house and offer condolence and aid, as is our custom miss emily met them at the door, dressed as usual and with no trace of grief on her face. she told them that her father was not dead. she did that for three days, with the ministers calling on her, and the doctors, trying to persuade her to let them dispose of the body. just as they were about to resort to law and force, she broke down, and they buried her father quickly. we did not say she was crazy then. we believed she had to do that. we remembered all the young men her father had driven away, and we knew that with nothing left, she would have to cling to that which had robbed her, as people will.  iii '
$ws.Range('AT14').Value = 'She was sick for a long time. when we saw her again, her hair was cut short, making her look like a girl, with a vague resemblance to those angels in colored church windows--sort of tragic and serene. the town had just let the contracts for paving the sidewalks, and in the summer after her father''s death they began the work. the construction company came with riggers and mules and machinery, and a foreman named homer barron, a yankee--a big, dark, ready man, with a big voice and eyes lighter than his face. the little boys would follow in groups to hear him cuss the riggers, and the riggers singing in time to the rise and fall of picks. pretty soon he knew everybody in town. whenever you heard a lot of laughing anywhere about the square, homer barron would be in the center of the group. '
$ws.Range('AT15').Value = 'Presently we began to see him and miss emily on sunday afternoons driving in the yellow-wheeled buggy and the matched team of bays from the livery stable. at first we were glad that miss emily would have an interest, because the ladies all said, "of course a grierson would not think seriously of a northerner, a day laborer." but there were still others, older people, who said that even grief could not cause a real lady to forget noblesse obligewithout calling it noblesse oblige. they just said, "poor emily. her kinsfolk should come to her." '
$ws.Range('AT16').Value = 'She had some kin in alabama; but years ago her father had fallen out with them over the estate of old lady wyatt, the crazy woman, and there was no communication between the two families. they had not even been represented at the funeral. '
$ws.Range('AT17').Value = 'And as soon as the old people said, "poor emily," the whispering began. "do you suppose it''s really so?" they said to one another. "of course it is. what else could . . ." this behind their hands; rustling of craned silk and satin behind jalousies closed upon the sun of sunday afternoon as the thin, swift clop-clop-clop of the matched team passed: "poor emily." she carried her head high enough--even when we believed that she was fallen. it was as if she demanded more than ever the recognition of her dignity as the last grierson; as if it had wanted that touch of earthiness to reaffirm her imperviousness. '
$ws.Range('AT18').Value = 'Like when she bought the rat poison, the arsenic. that was over a year after they had begun to say "poor emily," and while the two female cousins were visiting her. "i want some poison," she said to the druggist. she was over thirty then, still a slight woman, though thinner than usual, with cold, haughty black  eyes in a face the flesh of which was strained across the temples and about the eyesockets as you imagine a lighthouse-keeper''s face ought to look. "i want some poison," she said. "yes, miss emily. what kind? for rats and such? i''d recom--" "i want the best you have. i don''t care what kind." the druggist named several. "they''ll kill anything up to an elephant. but what you want is--" "arsenic," miss emily said. "is that a good one?" "is . . . arsenic? yes, ma''am. but what you want--" "i want arsenic." the druggist looked down at her. she looked back at him, erect, her face like a strained flag. "why, of course," the druggist said. "if that''s what you want. but the law requires you to tell what you are going to use it for." miss emily just stared at him, her head tilted back in order to look him eye for eye, until he looked away and went and got the arsenic and wrapped it up. the negro delivery boy brought her the package; the druggist didn''t come back. when she opened the package at home there was written on the box, under the skull and bones: "for rats."  iv '
$ws.Range('AT19').Value = 'So the next day we all said, "she will kill herself"; and we said it would be the best thing. '
$ws.Range('AT20').Value = 'When she had first begun to be seen with homer barron, we had said, "she will marry him." then we said, "she will persuade him yet," because homer himself had remarked--he liked men, and it was known that he drank with the younger men in the elks'' club-that he was not a marrying man. later we said, "poor emily" behind the jalousies as they passed on sunday afternoon in the glittering buggy, miss emily with her head high and homer barron with his hat cocked and a cigar in his teeth, reins and whip in a yellow glove. then some of the ladies began to say that it was a disgrace to the town and a bad example to the young people. the men did not want to interfere, '
$ws.Range('AT21').Value = 'but at last the ladies forced the baptist minister--miss emily''s people were episcopal-- to call upon her. he would never divulge what happened during that interview, but he refused to go back again. the  next sunday they again drove about the streets, '
$ws.Range('AT22').Value = 'and the following day the minister''s wife wrote to miss emily''s relations in alabama. '
$ws.Range('AT23').Value = 'So she had blood-kin under her roof again and we sat back to watch developments. at first nothing happened. then we were sure that they were to be married. '
$ws.Range('AT24').Value = 'We learned that miss emily had been to the jeweler''s and ordered a man''s toilet set in silver, with the letters h. b. on each piece. '
$ws.Range('AT25').Value = 'Two days later we learned that she had bought a complete outfit of men''s clothing, including a nightshirt, and we said, "they are married." we were really glad. we were glad because the two female cousins were even more grierson than miss emily had ever been. '
$ws.Range('AT26').Value = 'So we were not surprised when homer barron--the streets had been finished some time since--was gone. we were a little disappointed that there was not a public blowing-off, but we believed that he had gone on to prepare for miss emily''s coming, or to give her a chance to get rid of the cousins. (by that time it was a cabal, and we were all miss emily''s allies to help circumvent the cousins.) '
$ws.Range('AT27').Value = 'Sure enough, after another week they departed. '
$ws.Range('AT28').Value = 'And, as we had expected all along, within three days homer barron was back in town. a neighbor saw the negro man admit him at the kitchen door at dusk one evening. '
$ws.Range('AT29').Value = 'And that was the last we saw of homer barron. and of miss emily for some time. the negro man went in and out with the market basket, but the front door remained closed. now and then we would see her at a window for a moment, as the men did that night when they sprinkled the lime, but for almost six months she did not appear on the streets. then we knew that this was to be expected too; as if that quality of her father which had thwarted her woman''s life so many times had been too virulent and too furious to die. '
$ws.Range('AT30').Value = 'When we next saw miss emily, she had grown fat and her hair was turning gray. during the next few years it grew grayer and grayer until it attained an even pepper-and-salt iron-gray, when it ceased turning. up to the day of her death at seventy-four it was still that vigorous iron-gray, like the hair of an active man. '
$ws.Range('AT31').Value = 'From that time on her front door remained closed, save for a period of six or seven years, when she was about forty, during which she gave lessons in china-painting. she fitted up a studio in one of the downstairs rooms, where the daughters and granddaughters of colonel sartoris'' contemporaries were sent to her with the same regularity and in the same spirit that they were sent to church on sundays with a twenty-five-cent piece for the collection plate. meanwhile her taxes had been remitted.  '
$ws.Range('AT32').Value = 'Then the newer generation became the backbone and the spirit of the town, and the painting pupils grew up and fell away and did not send their children to her with boxes of color and tedious brushes and pictures cut from the ladies'' magazines. the front door closed upon the last one and remained closed for good. when the town got free postal delivery, miss emily alone refused to let them fasten the metal numbers above her door and attach a mailbox to it. she would not listen to them. daily, monthly, yearly we watched the negro grow grayer and more stooped, going in and out with the market basket. each december we sent her a tax notice, which would be returned by the post office a week later, unclaimed. now and then we would see her in one of the downstairs windows--she had evidently shut up the top floor of the house--like the carven torso of an idol in a niche, looking or not looking at us, we could never tell which. thus she passed from generation to generation--dear, inescapable, impervious, tranquil, and perverse. '
$ws.Range('AT33').Value = 'And so she died. fell ill in the house filled with dust and shadows, with only a doddering negro man to wait on her. we did not even know she was sick; we had long since given up trying to get any information from the negro he talked to no one, probably not even to her, for his voice had grown harsh and rusty, as if from disuse. she died in one of the downstairs rooms, in a heavy walnut bed with a curtain, her gray head propped on a pillow yellow and moldy with age and lack of sunlight.  v '
$ws.Range('AT34').Value = 'The negro met the first of the ladies at the front door and let them in, with their hushed, sibilant voices and their quick, curious glances, and then he disappeared. he walked right through the house and out the back and was not seen again. '
$ws.Range('AT35').Value = 'The two female cousins came at once. they held the funeral on the second day, with the town coming to look at miss emily beneath a mass of bought flowers, with the crayon face of her father musing profoundly above the bier and the ladies sibilant and macabre; and the very old men --some in their brushed confederate uniforms--on the porch and the lawn, talking of miss emily as if she had been a contemporary of theirs, believing that they had danced with her and courted her perhaps, confusing time with its mathematical progression, as the old do, to whom all the past is not a diminishing road but, instead, a huge meadow which  no winter ever quite touches, divided from them now by the narrow bottle-neck of the most recent decade of years. '
$ws.Range('AT36').Value = 'Already we knew that there was one room in that region above stairs which no one had seen in forty years, and which would have to be forced. they waited until miss emily was decently in the ground before they opened it. the violence of breaking down the door seemed to fill this room with pervading dust. a thin, acrid pall as of the tomb seemed to lie everywhere upon this room decked and furnished as for a bridal: upon the valance curtains of faded rose color, upon the rose-shaded lights, upon the dressing table, upon the delicate array of crystal and the man''s toilet things backed with tarnished silver, silver so tarnished that the monogram was obscured. among them lay a collar and tie, as if they had just been removed, which, lifted, left upon the surface a pale crescent in the dust. upon a chair hung the suit, carefully folded; beneath it the two mute shoes and the discarded socks. the man himself lay in the bed. for a long while we just stood there, looking down at the profound and fleshless grin. the body had apparently once lain in the attitude of an embrace, but now the long sleep that outlasts love, that conquers even the grimace of love, had cuckolded him. what was left of him, rotted beneath what was left of the nightshirt, had become inextricable from the bed in which he lay; and upon him and upon the pillow beside him lay that even coating of the patient and biding dust. then we noticed that in the second pillow was the indentation of a head. one of us lifted something from it, and leaning forward, that faint and invisible dust dry and acrid in the nostrils, we saw a long strand of irongray hair. '
